$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for Wins, Losses, Ties using same style as the neighboring header (AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record values for each data row (2-48)
$ws.Range("AD2:AD48").Value = 83
$ws.Range("AE2:AE48").Value = 79
$ws.Range("AF2:AF48").Value = 0

Write-Host "done"
